$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Insert the missing job-ID value for row 2 and renumber the following rows
# so the S.No sequence runs 1..9 starting at row 2 instead of row 3..10.
$ws.Range("A10").Value = 9
$ws.Range("A9").Value = 8
$ws.Range("A8").Value = 7
$ws.Range("A7").Value = 6
$ws.Range("A6").Value = 5
$ws.Range("A5").Value = 4
$ws.Range("A4").Value = 3
$ws.Range("A3").Value = 2
$ws.Range("A2").Value = 1

# Select A2:A10 (active cell A2) and scroll back so the view starts at A1
# (removing the previous C1 top-left / W2 selection).
$ws.Range("A2:A10").Select()
